# "optimize index case group": mark a large batch of index test cases as
# no longer testable (column B "Testable": y -> n), extend the sheet's
# AutoFilter / _FilterDatabase range to cover the full data extent, move
# the selection, and fix up the couple of cells whose highlight moved to
# a different row as a side effect of the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column B ("Testable"): y -> n for the affected rows -------------------
$ws.Range("B2:B21").Value = "n"
$ws.Range("B23:B28").Value = "n"
$ws.Range("B30:B43").Value = "n"
$ws.Range("B45").Value = "n"
$ws.Range("B47:B48").Value = "n"
$ws.Range("B50:B52").Value = "n"
$ws.Range("B57:B98").Value = "n"
$ws.Range("B100:B114").Value = "n"
$ws.Range("B119:B120").Value = "n"
$ws.Range("B179:B202").Value = "n"
$ws.Range("B204:B211").Value = "n"
$ws.Range("B213").Value = "n"
$ws.Range("B215:B216").Value = "n"
$ws.Range("B221:B222").Value = "n"
$ws.Range("B225:B227").Value = "n"
$ws.Range("B231:B255").Value = "n"

# --- Highlight moves from A119/A120 to A228 ---------------------------------
# Copy the still-yellow A119 formatting onto A228 first (before A119 loses it).
$ws.Cells.Item(119, 1).Copy()
$ws.Cells.Item(228, 1).PasteSpecial(-4122)

# Then clear the highlight on A119/A120, matching an already-unhighlighted cell.
$ws.Cells.Item(23, 1).Copy()
$ws.Cells.Item(119, 1).PasteSpecial(-4122)
$ws.Cells.Item(120, 1).PasteSpecial(-4122)

# --- Extend AutoFilter + the hidden _xlnm._FilterDatabase name over full data
$ws.AutoFilterMode = $false
$ws.Range("A1:L257").AutoFilter()
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$L`$257"
    }
}

# --- Move the active selection -------------------------------------------
$ws.Range("C261").Select()
